# Zaimplementowano barebone poprawnego zapisywania danych z openpyxl
# Write the rows that openpyxl's writer now emits: a data row (2), four
# blank-but-present rows (3-6), and three more data rows (7-9), across
# columns A:AU (47 columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(3,3,3,3,3,3,1,1,253,253,246,3,3,246,246,253,246,246,253,253,253,246,246,246,246,3,1,1,1,1,1,254,3,1,10,254,1,1,254,1,1,1,3,12,12,1,3)
$row7 = @(0,1,1,1,0,0,0,0,1,1,0,0,0,0,0,1,0,0,1,1,1,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0)
$row8 = @(49699,49192,49192,49192,32801,36897,32801,32801,0,0,32769,32769,32801,32769,32769,0,32769,32769,0,0,0,32769,32769,32769,32769,32801,32768,32769,32769,32769,32801,257,32801,32769,129,257,32769,49161,257,32769,32769,32769,32800,20617,4225,32769,32801)
$row9 = @(63,63,63,63,63,63,63,63,45,45,63,63,63,63,63,45,63,63,45,45,45,63,63,63,63,63,63,63,63,63,63,45,63,63,63,45,63,63,45,63,63,63,63,63,63,63,63)

for ($c = 1; $c -le $row2.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2[$c - 1]
}

# Rows 3-6 exist in the sheet (openpyxl wrote them out) but carry no cell
# data. Touching a row-level property that's already at its default value
# forces the row to be materialised without adding any cell content.
foreach ($r in 3..6) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

for ($c = 1; $c -le $row7.Length; $c++) {
    $ws.Cells.Item(7, $c).Value = $row7[$c - 1]
}
for ($c = 1; $c -le $row8.Length; $c++) {
    $ws.Cells.Item(8, $c).Value = $row8[$c - 1]
}
for ($c = 1; $c -le $row9.Length; $c++) {
    $ws.Cells.Item(9, $c).Value = $row9[$c - 1]
}
